$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.013.05'
$ws.Range("E2").Value = '  -1.58%  '
$ws.Range("D3").Value = '3.497.82'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '584.14'
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("D6").Value = '130.97'
$ws.Range("E6").Value = '  -2.70%  '
$ws.Range("D7").Value = '3.499.87'
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  -1.95%  '
$ws.Range("D10").Value = '0.124'
$ws.Range("E10").Value = '  -0.66%  '
$ws.Range("D11").Value = '7.12'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("D12").Value = '0.376'
$ws.Range("E12").Value = '  -3.09%  '
$ws.Range("D13").Value = '4.089.43'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '27.30'
$ws.Range("E14").Value = '  -1.53%  '
$ws.Range("D15").Value = '0.118'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '3.494.62'
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '0.0000177'
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").Value = '63.985.30'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("D19").Value = '9.75'
$ws.Range("E19").Value = '  -3.55%  '
$ws.Range("D20").Value = '14.08'
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("D21").Value = '5.61'
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D22").Value = '381.20'
$ws.Range("E22").Value = '  -3.02%  '
$ws.Range("D23").Value = '0.570'
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").Value = '3.635.73'
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").Value = '73.25'
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").Value = '0.0000115'
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("D28").Value = '1.55'
$ws.Range("E28").Value = '  -1.50%  '
$ws.Range("D29").Value = '7.44'
$ws.Range("E29").Value = '  -4.15%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").Value = '8.24'
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("D32").Value = '2.23'
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("D33").Value = '3.506.15'
$ws.Range("E33").Value = '  -0.65%  '
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = '23.32'
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("D36").Value = '0.144'
$ws.Range("E36").Value = '  -0.17%  '
$ws.Range("D37").Value = '5.29'
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = '1.55'
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("D39").Value = '6.86'
$ws.Range("E39").Value = '  -2.03%  '
$ws.Range("D40").Value = '159.33'
$ws.Range("E40").Value = '  -5.57%  '
$ws.Range("D41").Value = '0.0788'
$ws.Range("E41").Value = '  -3.47%  '
$ws.Range("D42").Value = '0.810'
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("D43").Value = '26.05'
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = '41.89'
$ws.Range("E45").Value = '  -2.50%  '
$ws.Range("D46").Value = '4.37'
$ws.Range("E46").Value = '  -1.49%  '
$ws.Range("D47").Value = '1.20'
$ws.Range("E47").Value = '  -5.77%  '
$ws.Range("D48").Value = '1.60'
$ws.Range("E48").Value = '  -3.04%  '
$ws.Range("D49").Value = '6.81'
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("D50").Value = '2.417.29'
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").Value = '0.897'
$ws.Range("E51").Value = '  -1.45%  '
